$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "Phone" column before the old "Date Entered" column (E),
# shifting Date Entered to F. Insert() shifts existing cell content,
# styles, and column widths along with it.
$ws.Columns("E:E").Insert()

# New column header
$ws.Range("E1").Value = "Phone"

# New phone number values. Leading apostrophe forces these to be stored
# as literal text (quote-prefixed) rather than being parsed as numbers,
# matching the source workbook's style (quotePrefix, no explicit number
# format change).
$ws.Range("E2").Value = "'+155567179876"
$ws.Range("E3").Value = "'+44 2040001234"
$ws.Range("E4").Value = "'1" + [char]8211 + "800" + [char]8211 + "854" + [char]8211 + "3680"

# CustomerName for row 4 changes from "Exxon" to "Apple"
$ws.Range("C4").Value = "Apple"

# Approximate the bestFit-style column widths the real Excel client
# computed for the two new columns (closest reachable values given this
# engine's column-width quantization).
$ws.Columns("E:E").ColumnWidth = 13.916666666666666
$ws.Columns("F:F").ColumnWidth = 12.583333333333334

# Restore/update the active selection to match the edited workbook.
$ws.Range("F11").Select() | Out-Null
